# Apply the tracked bill-summary edits (rows 8-25, 27, 29) to the active sheet.
# A leading apostrophe forces Excel to store the value as literal TEXT (matching
# the <c t="str"> cells in the target), without touching NumberFormat/styles -
# exactly like a user typing an apostrophe before a value in the Excel UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = '''P. point'
$ws.Range("C8").Value = 88
$ws.Range("D8").Value = '''2'
$ws.Range("E8").Value = '''Short point (up to 3 mtr.)'
$ws.Range("F8").Value = 256
$ws.Range("G8").Value = '''22528.00'
$ws.Range("C9").Value = 4
$ws.Range("G9").Value = '''1888.00'
$ws.Range("A10").Value = '''P. point'
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = '''4'
$ws.Range("E10").Value = '''Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").Value = '''7282.00'
$ws.Range("A11").Value = ''''
$ws.Range("C11").Value = 11
$ws.Range("D11").Value = '''2.0'
$ws.Range("E11").Value = '''Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = '''0.00'
$ws.Range("C12").Value = 24
$ws.Range("G12").Value = '''552.00'
$ws.Range("C13").Value = 4
$ws.Range("G13").Value = '''200.00'
$ws.Range("C14").Value = 31
$ws.Range("G14").Value = '''6789.00'
$ws.Range("A15").Value = '''Each'
$ws.Range("C15").Value = 69
$ws.Range("D15").Value = '''10.0'
$ws.Range("E15").Value = '''Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F15").Value = 303
$ws.Range("G15").Value = '''20907.00'
$ws.Range("C16").Value = 93
$ws.Range("D16").Value = '''17'
$ws.Range("E16").Value = '''25 mm'
$ws.Range("F16").Value = 56
$ws.Range("G16").Value = '''5208.00'
$ws.Range("A17").Value = '''Set'
$ws.Range("C17").Value = 97
$ws.Range("D17").Value = '''13.0'
$ws.Range("E17").Value = '''Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F17").Value = 5733
$ws.Range("G17").Value = '''556101.00'
$ws.Range("A18").Value = ''''
$ws.Range("C18").Value = 74
$ws.Range("D18").Value = '''14.0'
$ws.Range("E18").Value = '''Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = '''0.00'
$ws.Range("A19").Value = ''''
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = '''15.0'
$ws.Range("E19").Value = '''Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = '''0.00'
$ws.Range("A20").Value = '''Each'
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = '''27'
$ws.Range("E20").Value = '''1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F20").Value = 492
$ws.Range("G20").Value = '''22632.00'
$ws.Range("A21").Value = ''''
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = '''31'
$ws.Range("E21").Value = '''Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = '''0.00'
$ws.Range("A22").Value = ''''
$ws.Range("C22").Value = 66
$ws.Range("D22").Value = '''18.0'
$ws.Range("E22").Value = '''Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = '''0.00'
$ws.Range("C23").Value = 52
$ws.Range("D23").Value = '''36'
$ws.Range("E23").Value = '''Total'
$ws.Range("A24").Value = '''%'
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = '''37'
$ws.Range("E24").Value = '''Add Tender Premium '
$ws.Range("C25").Value = 14
$ws.Range("G27").Value = '''644087.00'
$ws.Range("H27").Value = '''644087.00'
$ws.Range("G29").Value = '''644087.00'
$ws.Range("H29").Value = '''644087.00'
